# Structures-Loads Table.xlsx update
# Rebuilds the loads table with mass/launch/deployed load columns, material,
# tensile strength figures, a revised support-structure row label, and a new
# reference block underneath (Falcon Heavy stage thrust figures, lunar
# gravity, assumed g's during launch).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Reference / notes block (rows 11-15) ---------------------------------
$ws.Range("A11").Value = "Falcon Heavy First Stage Thrust:"
$ws.Range("B11").Value = "7.6 MN"

$ws.Range("A12").Value = "Falcon Heavy Second Stage Thrust:"
$ws.Range("B12").Value = "934 kN"

# ---- Header row (row 2) ---------------------------------------------------
$ws.Range("B2").Value = "Mass (kg)"
$ws.Range("C2").Value = "Launch Loads (N)"
$ws.Range("D2").Value = "Deployed Loads (N)"
$ws.Range("E2").Value = "Material"
$ws.Range("F2").Value = "Tensile Strength (Pa)"
$ws.Range("A2:F2").Font.Bold = $true
$ws.Range("A2:F2").Font.Underline = $true

# ---- Row 4 label update -----------------------------------------------
$ws.Range("A4").Value = "Support Structure (legs,etc)"

# ---- Reference block continued --------------------------------------------
$ws.Range("A13").Value = "Lunar Gravity: "
$ws.Range("B13").Value = 1.62

$ws.Range("A14").Value = "Assumed gs during launch"
$ws.Range("B14").Value = 3

# ---- Row 3: Rigid Module Shell data ---------------------------------------
$ws.Range("B3").Value = 13000
$ws.Range("C3").Formula = '=$B$14*$B$15*B3'
$ws.Range("D3").Formula = '=B3*$B$13'
$ws.Range("E3").Value = "Aluminum"
$ws.Range("F3").Value = 276000000
$ws.Range("F3").NumberFormat = "0.00E+00"

# ---- Row 4: remaining formulas/material ------------------------------------
$ws.Range("C4:C8").Formula = '=$B$14*$B$15*B4'
$ws.Range("D4").Formula = '=D3/4'
$ws.Range("E4").Value = "Aluminum"
$ws.Range("F4").Value = 276000000
$ws.Range("F4").NumberFormat = "0.00E+00"

# ---- Row 8: Berthing Mechanism mass/loads ----------------------------------
$ws.Range("B8").Value = 1000
$ws.Range("B8").NumberFormat = "0.00E+00"
$ws.Range("D8").Value = 3900

# ---- Reference block final entry ------------------------------------------
$ws.Range("A15").Value = "Earth 1g"
$ws.Range("B15").Value = 9.81

# ---- Column widths ----------------------------------------------------
$ws.Range("A1:B1").ColumnWidth = 27.75
$ws.Range("C1").ColumnWidth = 22.75
$ws.Range("D1").ColumnWidth = 18.25
$ws.Range("E1").ColumnWidth = 12.584
$ws.Range("F1").ColumnWidth = 18.917

# ---- View / selection ------------------------------------------------------
$ws.Range("D5").Select()

# ---- Page setup -------------------------------------------------------
$ws.PageSetup.Orientation = 1
